$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Status text "Ready for handoff" -> "In Translation" on every sheet that
# shows it (Overview!E2/F2, zh-cn!C2, de-de!C2 all share the same string).
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# Narrow the "zh-cn"/"de-de" status columns on Overview (E, F) and the
# "Status" column (C) on the per-locale sheets from ~17.22 chars down to
# ~13.41 chars of stored width.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
